$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price-check snapshot was taken; its timestamp column is inserted right
# before the trailing "nom" / "url_produit" columns (shifting them one column
# to the right: GE->GF, GF->GG).
$ws.Range("GE:GE").EntireColumn.Insert()
$ws.Range("GE1").Value = "2026-02-05 15:32:19"

# For every product row that currently has a tracked price (column GD,
# formerly the latest snapshot), the new snapshot repeats the same price
# (nothing changed since the previous check). Rows without a price (GD empty)
# stay empty in the new column too - which is already the case right after
# the column insert, so nothing else needs to be done for them.
for ($r = 2; $r -le 80; $r++) {
    $price = $ws.Cells.Item($r, 186).Value()
    $ws.Cells.Item($r, 187).Value = $price
}
